$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Region -> tax rate mapping (derived from existing data)
$taxRateByRegion = @{
    "North" = 0.1
    "South" = 0.08
    "East"  = 0.09
    "West"  = 0.07000000000000001
}

# New invoice rows to append: Invoice_ID, Region, Sales_Amount
$newRows = @(
    @("INV006", "North", 2506),
    @("INV007", "North", 1600),
    @("INV008", "West", 1800),
    @("INV009", "West", 1900),
    @("INV0010", "East", 2000)
)

$startRow = 7
$r = $startRow

foreach ($row in $newRows) {
    $invoiceId = $row[0]
    $region = $row[1]
    $salesAmount = $row[2]
    $taxRate = $taxRateByRegion[$region]
    $taxAmount = $salesAmount * $taxRate
    $totalAmount = $salesAmount + $taxAmount

    $ws.Cells.Item($r, 1).Value = $invoiceId
    $ws.Cells.Item($r, 2).Value = $region
    $ws.Cells.Item($r, 3).Value = $salesAmount
    $ws.Cells.Item($r, 4).Value = $taxRate
    $ws.Cells.Item($r, 5).Value = $taxAmount
    $ws.Cells.Item($r, 6).Value = $totalAmount

    $r = $r + 1
}
